$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unlock just the cells we need to
# update so the rest of the protected-sheet XML (password hash, other
# protection flags) is left completely untouched.
$ws.Range("A11").Locked = $false
$ws.Range("D2:E8").Locked = $false

# Update the confidential notice date from 2021-04-26 to 2021-04-27
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.4886857271132336
$ws.Range("E2").Value = 0.002062422659150132

$ws.Range("D3").Value = 0.2516003380204045
$ws.Range("E3").Value = -0.00285266010554841

$ws.Range("D4").Value = 0.1002470633666092
$ws.Range("E4").Value = -0.0004798464491363452

$ws.Range("D5").Value = 0.1022002058893203
$ws.Range("E5").Value = 0.003455360478147096

$ws.Range("D6").Value = 0.02943457554750315
$ws.Range("E6").Value = 0.004766536964980528

$ws.Range("D7").Value = 0.02783209006292929
$ws.Range("E7").Value = 0.0006035458317612541

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.0007522805592263815

# Restore the original locked state now that the edits are applied.
$ws.Range("A11").Locked = $true
$ws.Range("D2:E8").Locked = $true
